# Generate Report for Handoff
#
# The localization pipeline has moved the job from "In Translation" to
# "Ready for handoff": update the Status text wherever it is reported
# (per-language tables + the Overview rollup) and refresh the associated
# timestamps to the moment the handoff report was (re)generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- zh-cn sheet: Status (C2) + Latest Handoff Datetime (H2) ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-25 11:01:27"

# --- de-de sheet: Status (C2) + Latest Handoff Datetime (H2) ---
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-25 11:01:33"

# --- Overview sheet: per-language status (E2/F2) + rollup generate date (G2) ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-25 11:01:33"

# The longer "Ready for handoff" text no longer fits the old Status column
# width, so the columns get widened (as Excel does on a content-driven
# resize) on every sheet that shows the Status column.
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
